$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 5

$ws.Cells.Item($row, 1).Value = 8
$ws.Cells.Item($row, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item($row, 3).Value = "Coquimbo"
$ws.Cells.Item($row, 4).Value = 44448
$ws.Cells.Item($row, 4).Style = $ws.Cells.Item(2, 4).Style
$ws.Cells.Item($row, 4).NumberFormat = $ws.Cells.Item(2, 4).NumberFormat
$ws.Cells.Item($row, 5).Value = 4
$ws.Cells.Item($row, 6).Value = "Fruta"
$ws.Cells.Item($row, 7).Value = 100101
$ws.Cells.Item($row, 8).Value = "Berries"
$ws.Cells.Item($row, 9).Value = 100101001
$ws.Cells.Item($row, 10).Value = "Arándano (blue)"
$ws.Cells.Item($row, 11).Value = "Sin especificar"
$ws.Cells.Item($row, 12).Value = "Primera"
$ws.Cells.Item($row, 13).Value = 100
$ws.Cells.Item($row, 14).Value = 14000
$ws.Cells.Item($row, 15).Value = 15000
$ws.Cells.Item($row, 16).Value = 14500
$ws.Cells.Item($row, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item($row, 18).Value = "Provincia de Limarí"
$ws.Cells.Item($row, 19).Value = 7250
$ws.Cells.Item($row, 20).Value = 2
